$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: insert new columns D,E,F (Corequisites, Concurrent, Recommended)
# and move "Terms Typically Offered" header to column G
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"
$ws.Range("G1").Value = "Terms Typically Offered"

# Rows 2-83: fill Corequisites(D)/Concurrent(E) with "NA", Recommended(F) with
# "NA" (or extracted "Recommended:" text split out of the old Prerequisites column),
# and shift the old "Terms Typically Offered" value (col D) into new column G.
$ws.Range("D2").Value = "NA"
$ws.Range("E2").Value = "NA"
$ws.Range("F2").Value = "NA"
$ws.Range("G2").Value = "F"
$ws.Range("D3").Value = "NA"
$ws.Range("E3").Value = "NA"
$ws.Range("F3").Value = "NA"
$ws.Range("G3").Value = "F, W, SP"
$ws.Range("D4").Value = "NA"
$ws.Range("E4").Value = "NA"
$ws.Range("F4").Value = "NA"
$ws.Range("G4").Value = "TBD"
$ws.Range("D5").Value = "NA"
$ws.Range("E5").Value = "NA"
$ws.Range("F5").Value = "NA"
$ws.Range("G5").Value = "TBD"
$ws.Range("D6").Value = "NA"
$ws.Range("E6").Value = "NA"
$ws.Range("F6").Value = "NA"
$ws.Range("G6").Value = "F, W, SP"
$ws.Range("D7").Value = "NA"
$ws.Range("E7").Value = "NA"
$ws.Range("F7").Value = "NA"
$ws.Range("G7").Value = "W"
$ws.Range("D8").Value = "NA"
$ws.Range("E8").Value = "NA"
$ws.Range("F8").Value = "NA"
$ws.Range("G8").Value = "F, W, SP"
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = "NA"
$ws.Range("F9").Value = "NA"
$ws.Range("G9").Value = "F, W, SP"
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = "NA"
$ws.Range("F10").Value = "NA"
$ws.Range("G10").Value = "F, SP"
$ws.Range("D11").Value = "NA"
$ws.Range("E11").Value = "NA"
$ws.Range("F11").Value = "NA"
$ws.Range("G11").Value = "F, SP"
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "NA"
$ws.Range("G12").Value = "W, SP"
$ws.Range("D13").Value = "NA"
$ws.Range("E13").Value = "NA"
$ws.Range("F13").Value = "NA"
$ws.Range("G13").Value = "F, W, SP"
$ws.Range("D14").Value = "NA"
$ws.Range("E14").Value = "NA"
$ws.Range("F14").Value = "NA"
$ws.Range("G14").Value = "F, W, SP"
$ws.Range("D15").Value = "NA"
$ws.Range("E15").Value = "NA"
$ws.Range("F15").Value = "NA"
$ws.Range("G15").Value = "SP"
$ws.Range("D16").Value = "NA"
$ws.Range("E16").Value = "NA"
$ws.Range("F16").Value = "NA"
$ws.Range("G16").Value = "F"
$ws.Range("D17").Value = "NA"
$ws.Range("E17").Value = "NA"
$ws.Range("F17").Value = "NA"
$ws.Range("G17").Value = "TBD"
$ws.Range("D18").Value = "NA"
$ws.Range("E18").Value = "NA"
$ws.Range("F18").Value = "NA"
$ws.Range("G18").Value = "F, W, SP"
$ws.Range("D19").Value = "NA"
$ws.Range("E19").Value = "NA"
$ws.Range("F19").Value = "NA"
$ws.Range("G19").Value = "W"
$ws.Range("D20").Value = "NA"
$ws.Range("E20").Value = "NA"
$ws.Range("F20").Value = "NA"
$ws.Range("G20").Value = "F, SP"
$ws.Range("D21").Value = "NA"
$ws.Range("E21").Value = "NA"
$ws.Range("F21").Value = "NA"
$ws.Range("G21").Value = "TBD"
$ws.Range("D22").Value = "NA"
$ws.Range("E22").Value = "NA"
$ws.Range("F22").Value = "NA"
$ws.Range("G22").Value = "F"
$ws.Range("D23").Value = "NA"
$ws.Range("E23").Value = "NA"
$ws.Range("F23").Value = "NA"
$ws.Range("G23").Value = "TBD"
$ws.Range("D24").Value = "NA"
$ws.Range("E24").Value = "NA"
$ws.Range("F24").Value = "NA"
$ws.Range("G24").Value = "F, W, SP"
$ws.Range("D25").Value = "NA"
$ws.Range("E25").Value = "NA"
$ws.Range("F25").Value = "NA"
$ws.Range("G25").Value = "W, SP"
$ws.Range("D26").Value = "NA"
$ws.Range("E26").Value = "NA"
$ws.Range("F26").Value = "NA"
$ws.Range("G26").Value = "F, W, SP"
$ws.Range("D27").Value = "NA"
$ws.Range("E27").Value = "NA"
$ws.Range("F27").Value = "NA"
$ws.Range("G27").Value = "W"
$ws.Range("D28").Value = "NA"
$ws.Range("E28").Value = "NA"
$ws.Range("F28").Value = "NA"
$ws.Range("G28").Value = "SP"
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("F29").Value = "NA"
$ws.Range("G29").Value = "F"
$ws.Range("C30").Value = "ASCI 229 and one of the CHEM 212, CHEM 216, CHEM 312, or CHEM 316."
$ws.Range("D30").Value = "NA"
$ws.Range("E30").Value = "NA"
$ws.Range("F30").Value = "NA"
$ws.Range("G30").Value = "W, SP"
$ws.Range("D31").Value = "NA"
$ws.Range("E31").Value = "NA"
$ws.Range("F31").Value = "NA"
$ws.Range("G31").Value = "W"
$ws.Range("D32").Value = "NA"
$ws.Range("E32").Value = "NA"
$ws.Range("F32").Value = "NA"
$ws.Range("G32").Value = "TBD"
$ws.Range("D33").Value = "NA"
$ws.Range("E33").Value = "NA"
$ws.Range("F33").Value = "NA"
$ws.Range("G33").Value = "SP"
$ws.Range("D34").Value = "NA"
$ws.Range("E34").Value = "NA"
$ws.Range("F34").Value = "NA"
$ws.Range("G34").Value = "F"
$ws.Range("D35").Value = "NA"
$ws.Range("E35").Value = "NA"
$ws.Range("F35").Value = "NA"
$ws.Range("G35").Value = "W"
$ws.Range("D36").Value = "NA"
$ws.Range("E36").Value = "NA"
$ws.Range("F36").Value = "NA"
$ws.Range("G36").Value = "W"
$ws.Range("D37").Value = "NA"
$ws.Range("E37").Value = "NA"
$ws.Range("F37").Value = "NA"
$ws.Range("G37").Value = "F, W, SP"
$ws.Range("D38").Value = "NA"
$ws.Range("E38").Value = "NA"
$ws.Range("F38").Value = "NA"
$ws.Range("G38").Value = "W"
$ws.Range("D39").Value = "NA"
$ws.Range("E39").Value = "NA"
$ws.Range("F39").Value = "NA"
$ws.Range("G39").Value = "W"
$ws.Range("D40").Value = "NA"
$ws.Range("E40").Value = "NA"
$ws.Range("F40").Value = "NA"
$ws.Range("G40").Value = "F, W, SP"
$ws.Range("D41").Value = "NA"
$ws.Range("E41").Value = "NA"
$ws.Range("F41").Value = "NA"
$ws.Range("G41").Value = "W"
$ws.Range("D42").Value = "NA"
$ws.Range("E42").Value = "NA"
$ws.Range("F42").Value = "NA"
$ws.Range("G42").Value = "SP"
$ws.Range("D43").Value = "NA"
$ws.Range("E43").Value = "NA"
$ws.Range("F43").Value = "NA"
$ws.Range("G43").Value = "SP"
$ws.Range("D44").Value = "NA"
$ws.Range("E44").Value = "NA"
$ws.Range("F44").Value = "NA"
$ws.Range("G44").Value = "W"
$ws.Range("D45").Value = "NA"
$ws.Range("E45").Value = "NA"
$ws.Range("F45").Value = "NA"
$ws.Range("G45").Value = "F, W, SP"
$ws.Range("D46").Value = "NA"
$ws.Range("E46").Value = "NA"
$ws.Range("F46").Value = "NA"
$ws.Range("G46").Value = "W"
$ws.Range("D47").Value = "NA"
$ws.Range("E47").Value = "NA"
$ws.Range("F47").Value = "NA"
$ws.Range("G47").Value = "F, W, SP"
$ws.Range("D48").Value = "NA"
$ws.Range("E48").Value = "NA"
$ws.Range("F48").Value = "NA"
$ws.Range("G48").Value = "F, W, SP"
$ws.Range("D49").Value = "NA"
$ws.Range("E49").Value = "NA"
$ws.Range("F49").Value = "NA"
$ws.Range("G49").Value = "W"
$ws.Range("D50").Value = "NA"
$ws.Range("E50").Value = "NA"
$ws.Range("F50").Value = "NA"
$ws.Range("G50").Value = "W"
$ws.Range("D51").Value = "NA"
$ws.Range("E51").Value = "NA"
$ws.Range("F51").Value = "NA"
$ws.Range("G51").Value = "TBD"
$ws.Range("D52").Value = "NA"
$ws.Range("E52").Value = "NA"
$ws.Range("F52").Value = "NA"
$ws.Range("G52").Value = "TBD"
$ws.Range("D53").Value = "NA"
$ws.Range("E53").Value = "NA"
$ws.Range("F53").Value = "NA"
$ws.Range("G53").Value = "F"
$ws.Range("D54").Value = "NA"
$ws.Range("E54").Value = "NA"
$ws.Range("F54").Value = "NA"
$ws.Range("G54").Value = "F"
$ws.Range("D55").Value = "NA"
$ws.Range("E55").Value = "NA"
$ws.Range("F55").Value = "NA"
$ws.Range("G55").Value = "F"
$ws.Range("C56").Value = "ASCI 229; ASCI 351; and ASCI 406."
$ws.Range("D56").Value = "NA"
$ws.Range("E56").Value = "NA"
$ws.Range("F56").Value = "ASCI 320 or CHEM 371; CHEM 327."
$ws.Range("G56").Value = "W "
$ws.Range("D57").Value = "NA"
$ws.Range("E57").Value = "NA"
$ws.Range("F57").Value = "NA"
$ws.Range("G57").Value = "SP"
$ws.Range("D58").Value = "NA"
$ws.Range("E58").Value = "NA"
$ws.Range("F58").Value = "NA"
$ws.Range("G58").Value = "W"
$ws.Range("D59").Value = "NA"
$ws.Range("E59").Value = "NA"
$ws.Range("F59").Value = "NA"
$ws.Range("G59").Value = "F, SP"
$ws.Range("D60").Value = "NA"
$ws.Range("E60").Value = "NA"
$ws.Range("F60").Value = "NA"
$ws.Range("G60").Value = "W, SP"
$ws.Range("D61").Value = "NA"
$ws.Range("E61").Value = "NA"
$ws.Range("F61").Value = "NA"
$ws.Range("G61").Value = "TBD"
$ws.Range("D62").Value = "NA"
$ws.Range("E62").Value = "NA"
$ws.Range("F62").Value = "NA"
$ws.Range("G62").Value = "W, SP"
$ws.Range("D63").Value = "NA"
$ws.Range("E63").Value = "NA"
$ws.Range("F63").Value = "NA"
$ws.Range("G63").Value = "F, SP"
$ws.Range("C64").Value = "ASCI 229."
$ws.Range("D64").Value = "NA"
$ws.Range("E64").Value = "NA"
$ws.Range("F64").Value = "ASCI 320, CHEM 371 or equivalent."
$ws.Range("G64").Value = "SP "
$ws.Range("D65").Value = "NA"
$ws.Range("E65").Value = "NA"
$ws.Range("F65").Value = "NA"
$ws.Range("G65").Value = "TBD"
$ws.Range("C66").Value = "ASCI 333 and ASCI 351."
$ws.Range("D66").Value = "NA"
$ws.Range("E66").Value = "NA"
$ws.Range("F66").Value = "ASCI 405 and ASCI 406."
$ws.Range("G66").Value = "TBD "
$ws.Range("D67").Value = "NA"
$ws.Range("E67").Value = "NA"
$ws.Range("F67").Value = "NA"
$ws.Range("G67").Value = "TBD"
$ws.Range("D68").Value = "NA"
$ws.Range("E68").Value = "NA"
$ws.Range("F68").Value = "NA"
$ws.Range("G68").Value = "SU"
$ws.Range("D69").Value = "NA"
$ws.Range("E69").Value = "NA"
$ws.Range("F69").Value = "NA"
$ws.Range("G69").Value = "TBD"
$ws.Range("D70").Value = "NA"
$ws.Range("E70").Value = "NA"
$ws.Range("F70").Value = "NA"
$ws.Range("G70").Value = "TBD"
$ws.Range("C71").Value = "Senior standing, ASCI 363 and consent of instructor."
$ws.Range("D71").Value = "NA"
$ws.Range("E71").Value = "NA"
$ws.Range("F71").Value = "one course in statistics."
$ws.Range("G71").Value = "F, W, SP "
$ws.Range("D72").Value = "NA"
$ws.Range("E72").Value = "NA"
$ws.Range("F72").Value = "NA"
$ws.Range("G72").Value = "F, W, SP"
$ws.Range("D73").Value = "NA"
$ws.Range("E73").Value = "NA"
$ws.Range("F73").Value = "NA"
$ws.Range("G73").Value = "F, SP"
$ws.Range("D74").Value = "NA"
$ws.Range("E74").Value = "NA"
$ws.Range("F74").Value = "NA"
$ws.Range("G74").Value = "F, SP"
$ws.Range("D75").Value = "NA"
$ws.Range("E75").Value = "NA"
$ws.Range("F75").Value = "NA"
$ws.Range("G75").Value = "F, W, SP"
$ws.Range("D76").Value = "NA"
$ws.Range("E76").Value = "NA"
$ws.Range("F76").Value = "NA"
$ws.Range("G76").Value = "F, W, SP"
$ws.Range("C77").Value = "ASCI 320, or CHEM 313 or CHEM 371, and one of the ASCI 346, or ASCI 350, or ASCI 355, or DSCI 301, or consent of instructor."
$ws.Range("D77").Value = "NA"
$ws.Range("E77").Value = "NA"
$ws.Range("F77").Value = "NA"
$ws.Range("G77").Value = "TBD"
$ws.Range("D78").Value = "NA"
$ws.Range("E78").Value = "NA"
$ws.Range("F78").Value = "NA"
$ws.Range("G78").Value = "TBD"
$ws.Range("D79").Value = "NA"
$ws.Range("E79").Value = "NA"
$ws.Range("F79").Value = "NA"
$ws.Range("G79").Value = "TBD"
$ws.Range("D80").Value = "NA"
$ws.Range("E80").Value = "NA"
$ws.Range("F80").Value = "NA"
$ws.Range("G80").Value = "TBD"
$ws.Range("D81").Value = "NA"
$ws.Range("E81").Value = "NA"
$ws.Range("F81").Value = "NA"
$ws.Range("G81").Value = "F, W, SP"
$ws.Range("D82").Value = "NA"
$ws.Range("E82").Value = "NA"
$ws.Range("F82").Value = "NA"
$ws.Range("G82").Value = "TBD"
$ws.Range("D83").Value = "NA"
$ws.Range("E83").Value = "NA"
$ws.Range("F83").Value = "NA"
$ws.Range("G83").Value = "TBD"
